$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Genotyp paragraph
#    "Genotyp* (nur für HCV)" -> "Genotyp*" (keep italic run)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Genotyp* (nur für HCV)", $false, $false, $false, $false, $false, $true, 1, $false, "Genotyp*", 2) | Out-Null

#    insert a new italic run containing the opening-quote text that used
#    to live at the front of the following (non-italic) run
$r = $d.Content
$r.Find.Execute("Genotyp*“ wähle den entsprechenden Genotyp aus dem Drop-down Menü aus. Falls") | Out-Null
$r2 = $d.Content
$r2.Find.Execute("Genotyp*") | Out-Null
$r2.Collapse(0)
$r2.InsertAfter([string][char]0x201C + " ")
$r2.Font.Italic = $true

#    drop the now-duplicated opening quote + space from the following run
$d.Content.Find.Execute("“ wähle den entsprechenden Genotyp aus dem", $false, $false, $false, $false, $false, $true, 1, $false, "wähle den entsprechenden Genotyp aus dem", 2) | Out-Null

# ---------------------------------------------------------------------
#    " Menü aus. Falls es sich um einen anderen Virus als HCV handelt
#    einfach „" -> " Menü aus. Falls dieser nicht bekannt ist, einfach „"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(" Menü aus. Falls es sich um einen anderen Virus als HCV handelt einfach „", $false, $false, $false, $false, $false, $true, 1, $false, " Menü aus. Falls dieser nicht bekannt ist, einfach „", 2) | Out-Null

# ---------------------------------------------------------------------
#    "“ auswählen (dieses Feld nicht leer lassen)." -> "“ auswählen (Feld nicht leer lassen)."
# ---------------------------------------------------------------------
$d.Content.Find.Execute("“ auswählen (dieses Feld nicht leer lassen).", $false, $false, $false, $false, $false, $true, 1, $false, "“ auswählen (Feld nicht leer lassen).", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Target paragraph
#    "Target*" -> "Target" (drop the asterisk, still italic)
#    relocate the "_GoBack" bookmark to sit right after "Target"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Target*“ die amplifizierte Region aus dem Drop-down Menü auswählen.", $false, $false, $false, $false, $false, $true, 1, $false, "Target“ die amplifizierte Region aus dem Drop-down Menü auswählen.", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Worksheet paragraph - merge the two runs ("Im Workshee" + "t „") into
#    a single run, removing the "_GoBack" bookmark that used to sit
#    between them.
# ---------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()
$d.Content.Find.Execute(" Im Workshee" + [string][char]0x00A4, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# re-add the bookmark right after "Target" in the Target paragraph
$tr = $d.Content
$tr.Find.Execute("Target") | Out-Null
$tr.Collapse(0)
$d.Bookmarks.Add("_GoBack", $tr) | Out-Null
